# Fixed issue #13 Permitir que en los ficheros de metadatos dos columnas se puedan
# relacionar para crear SKOS jerarquicos.
#
# A new metadata row is inserted right after the header row (row 1). This new
# row 2 holds the "bare" (unqualified) slug identifier for each column -
# i.e. the same identifier used in row 3 (previously row 2) but without its
# "iaest-measure:" / "sdmx-dimension:" namespace prefix - so that two columns
# can reference each other by these plain names to build SKOS hierarchies.
#
# The previous extra annotation row that only carried "mapping-ano.xlsx" in
# column U is dropped entirely; column U's row 5 (the datatype row) now simply
# reads "xsd:date" like the rest of that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the old row 2 (old rows 2-5 shift down to 3-6).
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the bare slug identifiers.
$ws.Range("A2").Value = "desempleo-perceptores"
$ws.Range("B2").Value = "pension-media-por-percepcion"
$ws.Range("C2").Value = "salario-perceptores"
$ws.Range("D2").Value = "salario-medio-por-percepcion"
$ws.Range("E2").Value = "desempleo-medio-por-percepcion"
$ws.Range("F2").Value = "pension-percepciones"
$ws.Range("G2").Value = "provincia-nombre"
$ws.Range("H2").Value = "pension-media-por-persona"
$ws.Range("I2").Value = "municipio-nombre"
$ws.Range("J2").Value = "salario-retribucion"
$ws.Range("K2").Value = "salario-medio-anual"
$ws.Range("L2").Value = "salario-percepciones-por-persona"
$ws.Range("M2").Value = "pension-retribucion"
$ws.Range("N2").Value = "desempleo-retribucion"
$ws.Range("O2").Value = "desempleo-prestacion-media-anual"
$ws.Range("P2").Value = "salario-percepciones"
$ws.Range("Q2").Value = "salario-medio-por-persona"
$ws.Range("R2").Value = "provincia-codigo"
$ws.Range("S2").Value = "municipio-codigo"
$ws.Range("T2").Value = "pension-percepciones-por-persona"
$ws.Range("U2").Value = "ano"
$ws.Range("V2").Value = "pension-perceptores"
$ws.Range("W2").Value = "ccaa"

# The old trailing row (now shifted to row 6) only held the obsolete
# "mapping-ano.xlsx" marker in column U; delete that row entirely so the
# sheet ends at row 5 again (now matching the shifted datatype row).
$ws.Rows.Item(6).Delete()
